# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 356f4a7e-... row across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 356f4a7e-... (row 3)
$wsOverview.Range("G3").Value = "2016-11-03 19:14:29"

# zh-cn sheet: Correspond Handoff / Handback Datetime for 356f4a7e-... (row 3)
$wsZhCn.Range("H3").Value = "2016-11-03 19:14:16"
$wsZhCn.Range("K3").Value = "2016-11-03 19:15:05"

# de-de sheet: Correspond Handback DateTime for 356f4a7e-... (row 3)
$wsDeDe.Range("K3").Value = "2016-11-03 19:15:23"
